# Updates the crypto price/volume table with refreshed values scraped on
# Thu Feb  9 21:51:36 UTC 2023. Only the "Price" (D) and "Volume(1h)" (E)
# columns change; everything else (coin name, link, date, hour) stays put.
# The source sheet stores these as plain text (e.g. "308.96", "-5.71%"),
# not numbers, so each write is forced to stay text (NumberFormat "@" while
# assigning, then restored to the default "Normal" style) instead of being
# auto-converted to a number/percentage by Excel's input parsing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    "D2" = "308.07"
    "E2" = "-5.99%"
    "D3" = "40.46"
    "E3" = "-8.70%"
    "D4" = "5.061"
    "E4" = "-4.71%"
    "D5" = "0.07789"
    "E5" = "-6.96%"
    "D6" = "4.330"
    "E6" = "-1.63%"
    "D7" = "1.658"
    "E7" = "-13.76%"
    "D8" = "0.9128"
    "E8" = "-6.04%"
    "D9" = "0.1033"
    "E9" = "-8.92%"
    "D10" = "0.1748"
    "E10" = "-7.98%"
    "D11" = "0.08978"
    "E11" = "-6.97%"
    "D12" = "0.04451"
    "E12" = "-3.28%"
    "D13" = "7.111"
    "E13" = "-16.61%"
    "D14" = "0.1061"
    "E14" = "-0.04%"
    "D15" = "0.001252"
    "E15" = "-3.64%"
    "D16" = "0.006010"
    "E16" = "4.91%"
    "D17" = "3.364"
    "E17" = "-0.68%"
    "E18" = "0.58%"
    "D19" = "0.3366"
    "E19" = "0.27%"
    "D20" = "0.1384"
    "E20" = "-0.24%"
    "D21" = "0.2862"
    "E21" = "11.24%"
    "D22" = "0.04174"
    "E22" = "0.54%"
    "D23" = "0.001218"
    "E23" = "-1.19%"
    "D24" = "0.004100"
    "E24" = "-7.60%"
    "D25" = "0.0001231"
    "E25" = "-5.39%"
    "D26" = "0.0003005"
    "E26" = "0.90%"
    "D38" = "0.02397"
    "E38" = "-11.76%"
    "D39" = "0.05220"
    "E39" = "-6.86%"
    "D40" = "0.007989"
    "E40" = "1.87%"
    "D41" = "0.1333"
    "E41" = "-5.61%"
    "D42" = "0.007509"
    "E42" = "2.80%"
    "D43" = "0.002022"
    "E43" = "-1.46%"
    "D44" = "0.008094"
    "E44" = "-6.88%"
    "D45" = "0.3340"
    "E45" = "-4.85%"
    "D46" = "0.00006752"
    "E46" = "-2.12%"
    "D47" = "0.00000000757"
    "E47" = "0.92%"
    "D48" = "0.003342"
    "E48" = "-4.37%"
    "D49" = "0.004136"
    "E49" = "17.20%"
    "D50" = "0.00002118"
    "E50" = "0.92%"
    "D51" = "0.0002018"
    "E51" = "0.92%"
}

foreach ($key in $updates.Keys) {
    $cell = $ws.Range($key)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$key]
    $cell.Style = "Normal"
}
